$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4: DDR4 DIMM no. 16 data, DDR4 4x4 Spatial config for banks
$ws.Range("A4").Value = "05_24_2024"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Coarse"
$ws.Range("D4").Value = 15

# Match the selection left after the edit
$ws.Range("D4").Select()
